$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) B12:B16 ("Barkod") were saved as TEXT with a leading zero, e.g. "09225074".
#    Correct them to plain numbers (no leading zero), matching B2:B11 above.
# ---------------------------------------------------------------------------
$ws.Cells.Item(12, 2).Value = 9225074
$ws.Cells.Item(13, 2).Value = 9218978
$ws.Cells.Item(14, 2).Value = 9216193
$ws.Cells.Item(15, 2).Value = 9188047
$ws.Cells.Item(16, 2).Value = 9175143

# ---------------------------------------------------------------------------
# 2) Append 5 new product rows (17-21) after the existing data (O16 was the
#    previous last cell; sheet now spans A1:O21).
#    A leading apostrophe forces TEXT, which keeps the zero-padded barcode in B
#    as text (instead of Excel auto-detecting it as a number) and keeps the
#    untouched image-url cells as empty TEXT cells (instead of fully blank ones),
#    matching how the rest of the sheet is authored.
# ---------------------------------------------------------------------------
$newRows = @(
    [ordered]@{ ProductID=''; Barkod='09225074'; ProductName='XMART TAŞINABİLİR MULTİMEDYA PROJEKTÖR MPP-40'; Price=169; Currency='BGN'; Category='Projektörler'; Brand='XMART'; AnaGorsel='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09225074/67188f8bd2191.jpg.webp'; Image1='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09225074/67188f8d5e802.jpg.webp'; Image2='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09225074/67188f8f76bb7.jpg.webp'; Image3='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09225074/67188f9188007.jpg.webp'; Image4=''; Image5=''; DigerGorseller=''; ProductURL='https://www.technomarket.bg/proektori/xmart-portable-multimedia-projector-mpp-40-09225074' }
    [ordered]@{ ProductID=''; Barkod='09218978'; ProductName='EPSON CO-W01'; Price=699; Currency='BGN'; Category='Projektörler'; Brand='EPSON'; AnaGorsel='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09218978/6618e8c99b73e.jpg.webp'; Image1='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09218978/6618e8cbaed42.jpg.webp'; Image2='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09218978/6618e8d0c2b08.jpg.webp'; Image3='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09218978/6618e8d32d971.jpg.webp'; Image4='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09218978/6618e8d59627a.jpg.webp'; Image5='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09218978/6618e8d8b2bb5.jpg.webp'; DigerGorseller=''; ProductURL='https://www.technomarket.bg/proektori/epson-co-w01-09218978' }
    [ordered]@{ ProductID=''; Barkod='09216193'; ProductName='SAMSUNG SP-LFF3CLA SERBEST STİL AKILLI'; Price=1099; Currency='BGN'; Category='Projektörler'; Brand='SAMSUNG'; AnaGorsel='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09216193/64eceb152806d.jpg.webp'; Image1='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09216193/64eceb16b6442.jpg.webp'; Image2='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09216193/64eceb1869b71.jpg.webp'; Image3='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09216193/64eceb1a0cf19.jpg.webp'; Image4='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09216193/64eceb1d67ef1.jpg.webp'; Image5='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09216193/64eceb20b69a1.jpg.webp'; DigerGorseller=''; ProductURL='https://www.technomarket.bg/proektori/samsung-sp-lff3cla-the-freestyle-smart-09216193' }
    [ordered]@{ ProductID=''; Barkod='09188047'; ProductName='ACER PROJEKTÖR X128HP'; Price=819; Currency='BGN'; Category='Projektörler'; Brand='ACER'; AnaGorsel='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09188047/6077de34291d4.jpg.webp'; Image1=''; Image2=''; Image3=''; Image4=''; Image5=''; DigerGorseller=''; ProductURL='https://www.technomarket.bg/proektori/acer-projector-x128hp-09188047' }
    [ordered]@{ ProductID=''; Barkod='09175143'; ProductName='ACER PROJEKTÖR X1126AH'; Price=728; Currency='BGN'; Category='Projektörler'; Brand='ACER'; AnaGorsel='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09175143/5d415857e49d5.jpg.webp'; Image1='https://cdn.technomarket.bg/ng/media/cache/mid_thumb/uploads/library/product/09175143/5d4158592eb69.jpg.webp'; Image2=''; Image3=''; Image4=''; Image5=''; DigerGorseller=''; ProductURL='https://www.technomarket.bg/proektori/acer-projector-x1126ah-09175143' }
)

# Columns A..O in sheet order
$keys = @("ProductID","Barkod","ProductName","Price","Currency","Category","Brand","AnaGorsel","Image1","Image2","Image3","Image4","Image5","DigerGorseller","ProductURL")

$r = 17
foreach ($row in $newRows) {
    for ($i = 0; $i -lt $keys.Length; $i++) {
        $key = $keys[$i]
        $val = $row[$key]
        if ($val -is [string]) {
            # leading apostrophe -> force text, even for empty / numeric-looking strings
            $ws.Cells.Item($r, $i + 1).Value = "'" + $val
        } else {
            $ws.Cells.Item($r, $i + 1).Value = $val
        }
    }
    $r++
}
